$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "correct" column (I) with header and per-row answer key (a/b/c/d)
# matching which of columns E/F/G/H (a/b/c/d) holds the correct response.
$ws.Range("I1").Value = "correct"
$ws.Range("I2").Value = "d"
$ws.Range("I3").Value = "c"
$ws.Range("I4").Value = "c"
$ws.Range("I5").Value = "a"
$ws.Range("I6").Value = "a"
$ws.Range("I7").Value = "c"
$ws.Range("I8").Value = "b"
$ws.Range("I9").Value = "a"
$ws.Range("I10").Value = "b"
$ws.Range("I11").Value = "a"
$ws.Range("I12").Value = "c"
$ws.Range("I13").Value = "b"
$ws.Range("I14").Value = "d"
$ws.Range("I15").Value = "b"
$ws.Range("I16").Value = "c"
$ws.Range("I17").Value = "a"
$ws.Range("I18").Value = "b"
$ws.Range("I19").Value = "c"
$ws.Range("I20").Value = "a"
$ws.Range("I21").Value = "c"
$ws.Range("I22").Value = "a"
$ws.Range("I23").Value = "a"
$ws.Range("I24").Value = "d"
$ws.Range("I25").Value = "a"
$ws.Range("I26").Value = "c"
$ws.Range("I27").Value = "d"
$ws.Range("I28").Value = "b"
$ws.Range("I29").Value = "a"

# Match the text-formatted style used by the other response columns (E:H)
$ws.Range("I1:I29").NumberFormat = "@"

# Update the selected cell to match the author's final selection
[void]$ws.Range("B5").Select()
